$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Relabel the QA report keys (rows shifted: a new "correct" key was
# inserted ahead of the existing "title"/"msg" keys).
$ws.Range("A2").Value = "qa_report_result_correct"
$ws.Range("A3").Value = "qa_report_result_title"
$ws.Range("A6").Value = "qa_report_correct"
$ws.Range("A7").Value = "qa_report_result_msg"

# qa_result flag flipped from TRUE to FALSE. Prefix with an apostrophe so
# Excel stores it as literal text (matching the existing cell), not a
# boolean.
$ws.Range("B13").Value = "'FALSE"
